# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Re-sorts / refreshes the "Estado de Cuenta" detail rows (B16:G24) on
# Hoja1 with the updated worker roster: a new worker (CINDY PATRICIA DIAZ
# OTERO, CC 1069489588) is added with four overdue periods (1701-1704),
# and the existing workers' period/valor-mora figures are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: CINDY PATRICIA DIAZ OTERO - periodo 1701
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1069489588"
$ws.Range("D16").Value = "CINDY PATRICIA DIAZ OTERO"
$ws.Range("E16").Value = "1701"
$ws.Range("F16").Value = 31477
$ws.Range("G16").Value = 1475500

# Row 17: CINDY PATRICIA DIAZ OTERO - periodo 1702
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1069489588"
$ws.Range("D17").Value = "CINDY PATRICIA DIAZ OTERO"
$ws.Range("E17").Value = "1702"
$ws.Range("F17").Value = 59020
$ws.Range("G17").Value = 1475500

# Row 18: DALGY MEJIA GUERRERO - periodo 1703
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "22801335"
$ws.Range("D18").Value = "DALGY MEJIA GUERRERO"
$ws.Range("E18").Value = "1703"
$ws.Range("F18").Value = 3935
$ws.Range("G18").Value = 737717

# Row 19: CINDY PATRICIA DIAZ OTERO - periodo 1703
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1069489588"
$ws.Range("D19").Value = "CINDY PATRICIA DIAZ OTERO"
$ws.Range("E19").Value = "1703"
$ws.Range("F19").Value = 59020
$ws.Range("G19").Value = 1475500

# Row 20: KATTY MILENA SEÑA GARCIA - periodo 1703
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143397062"
$ws.Range("D20").Value = "KATTY MILENA SEÑA GARCIA"
$ws.Range("E20").Value = "1703"
$ws.Range("F20").Value = 21640
$ws.Range("G20").Value = 737717

# Row 21: DALGY MEJIA GUERRERO - periodo 1704
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "22801335"
$ws.Range("D21").Value = "DALGY MEJIA GUERRERO"
$ws.Range("E21").Value = "1704"
$ws.Range("F21").Value = 29509
$ws.Range("G21").Value = 737717

# Row 22: XIOMARA HERAZO RODRIGUEZ - periodo 1704
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "45565136"
$ws.Range("D22").Value = "XIOMARA HERAZO RODRIGUEZ"
$ws.Range("E22").Value = "1704"
$ws.Range("F22").Value = 44143
$ws.Range("G22").Value = 1103575

# Row 23: CINDY PATRICIA DIAZ OTERO - periodo 1704
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1069489588"
$ws.Range("D23").Value = "CINDY PATRICIA DIAZ OTERO"
$ws.Range("E23").Value = "1704"
$ws.Range("F23").Value = 59020
$ws.Range("G23").Value = 1475500

# Row 24: KATTY MILENA SEÑA GARCIA - periodo 1704
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1143397062"
$ws.Range("D24").Value = "KATTY MILENA SEÑA GARCIA"
$ws.Range("E24").Value = "1704"
$ws.Range("F24").Value = 29509
$ws.Range("G24").Value = 737717
